# aggiornamento 15, 16, 17 marzo
# Append three new daily rows (227-229) to the bottom of the time series,
# continuing directly after the last existing row (226).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data to append: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @{ Row = 227; Date = 44301; B = 0; C = 5; D = 417.0141784820684 },
    @{ Row = 228; Date = 44302; B = 0; C = 2; D = 166.8056713928273 },
    @{ Row = 229; Date = 44303; B = 0; C = 2; D = 166.8056713928273 }
)

# Column A carries a special date/time style (border + centered + date format).
# Copy the formatting from the last existing row so the new cells match it,
# then fill in the values.
$lastRow = 226

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$lastRow").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

    $ws.Range("A$row").Value = $r.Date
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
}

$excel.CutCopyMode = $false
